$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 42 — this pushes the existing rows 42..95 down to 43..96
# (matching the diff, which shows every data row from the old row 42 onward
# shifting down by one position) and grows the used range to A1:R96.
$ws.Rows(42).Insert()

# Populate the newly inserted row 42 with its data.
$ws.Range("A42").Value = 7
$ws.Range("B42").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C42").Value = "Ñuble"
$ws.Range("D42").Value = 44494
$ws.Range("E42").Value = 16
$ws.Range("F42").Value = 100112045
$ws.Range("G42").Value = "Zapallo"
$ws.Range("H42").Value = "Camote"
$ws.Range("I42").Value = "1a (guarda)"
$ws.Range("J42").Value = 200
$ws.Range("K42").Value = 800
$ws.Range("L42").Value = 900
$ws.Range("M42").Value = 850
$ws.Range("N42").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O42").Value = "Región del Maule"
$ws.Range("P42").Value = 850
$ws.Range("Q42").Value = 1
$ws.Range("R42").Value = "Hortaliza"
